# prospects import example
#
# The "Phone" entry for row 3 was re-keyed from a quoted literal
# ("+60 012 9892525") to a plain, space-prefixed value that is kept as
# text (so Excel doesn't try to reinterpret/format it as a number), and
# the active selection moved from D17 to E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rngPhone = $ws.Range("E3")
$rngPhone.NumberFormat = "@"
$rngPhone.Value = " +60 012 9892525"

$ws.Range("E11").Select()
